$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $addr, $val)
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue $ws "D2" "293.12"
Set-TextValue $ws "E2" "-0.37%"
Set-TextValue $ws "G2" "11"

Set-TextValue $ws "D3" "40.34"
Set-TextValue $ws "E3" "0.83%"
Set-TextValue $ws "G3" "11"

Set-TextValue $ws "D4" "5.002"
Set-TextValue $ws "G4" "11"

Set-TextValue $ws "D5" "0.07334"
Set-TextValue $ws "E5" "-0.68%"
Set-TextValue $ws "G5" "11"

$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue $ws "D6" "1.560"
Set-TextValue $ws "E6" "1.65%"
Set-TextValue $ws "G6" "11"

$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws "D7" "0.9232"
Set-TextValue $ws "E7" "0.15%"
Set-TextValue $ws "G7" "11"

$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue $ws "D8" "2.345"
Set-TextValue $ws "E8" "-2.25%"
Set-TextValue $ws "G8" "11"

$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue $ws "D9" "0.1180"
Set-TextValue $ws "E9" "0.83%"
Set-TextValue $ws "G9" "11"

$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws "D10" "0.1809"
Set-TextValue $ws "E10" "2.98%"
Set-TextValue $ws "G10" "11"

$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws "D11" "0.04389"
Set-TextValue $ws "E11" "5.15%"
Set-TextValue $ws "G11" "11"

$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws "D12" "0.08795"
Set-TextValue $ws "E12" "1.65%"
Set-TextValue $ws "G12" "11"

$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws "D13" "0.1054"
Set-TextValue $ws "E13" "0.10%"
Set-TextValue $ws "G13" "11"

$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws "D14" "0.001272"
Set-TextValue $ws "E14" "1.53%"
Set-TextValue $ws "G14" "11"

$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws "D15" "0.005859"
Set-TextValue $ws "E15" "1.04%"
Set-TextValue $ws "G15" "11"

$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws "D16" "3.344"
Set-TextValue $ws "E16" "-0.86%"
Set-TextValue $ws "G16" "11"

$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue $ws "D17" "4.285"
Set-TextValue $ws "E17" "-0.85%"
Set-TextValue $ws "G17" "11"

Set-TextValue $ws "D18" "0.3306"
Set-TextValue $ws "E18" "0.26%"
Set-TextValue $ws "G18" "11"

Set-TextValue $ws "D19" "7.875"
Set-TextValue $ws "E19" "3.36%"
Set-TextValue $ws "G19" "11"

Set-TextValue $ws "E20" "2.52%"
Set-TextValue $ws "G20" "11"

Set-TextValue $ws "D21" "0.2800"
Set-TextValue $ws "E21" "-0.64%"
Set-TextValue $ws "G21" "11"

Set-TextValue $ws "D22" "0.03923"
Set-TextValue $ws "E22" "2.45%"
Set-TextValue $ws "G22" "11"

Set-TextValue $ws "D23" "0.001262"
Set-TextValue $ws "E23" "-1.84%"
Set-TextValue $ws "G23" "11"

Set-TextValue $ws "E24" "1.14%"
Set-TextValue $ws "G24" "11"

Set-TextValue $ws "E25" "-8.07%"
Set-TextValue $ws "G25" "11"

Set-TextValue $ws "D26" "0.0003723"
Set-TextValue $ws "E26" "-0.39%"
Set-TextValue $ws "G26" "11"

Set-TextValue $ws "G27" "11"

Set-TextValue $ws "G28" "11"

Set-TextValue $ws "G29" "11"

Set-TextValue $ws "G30" "11"

Set-TextValue $ws "G31" "11"

Set-TextValue $ws "G32" "11"

Set-TextValue $ws "G33" "11"

Set-TextValue $ws "G34" "11"

Set-TextValue $ws "G35" "11"

Set-TextValue $ws "G36" "11"

Set-TextValue $ws "G37" "11"

Set-TextValue $ws "D38" "0.02340"
Set-TextValue $ws "E38" "1.18%"
Set-TextValue $ws "G38" "11"

Set-TextValue $ws "D39" "0.05099"
Set-TextValue $ws "E39" "2.16%"
Set-TextValue $ws "G39" "11"

Set-TextValue $ws "D40" "0.005932"
Set-TextValue $ws "E40" "47.71%"
Set-TextValue $ws "G40" "11"

Set-TextValue $ws "D41" "0.007861"
Set-TextValue $ws "E41" "1.65%"
Set-TextValue $ws "G41" "11"

Set-TextValue $ws "D42" "0.1289"
Set-TextValue $ws "E42" "1.33%"
Set-TextValue $ws "G42" "11"

Set-TextValue $ws "G43" "11"

Set-TextValue $ws "D44" "0.008025"
Set-TextValue $ws "E44" "1.76%"
Set-TextValue $ws "G44" "11"

Set-TextValue $ws "D45" "0.2916"
Set-TextValue $ws "E45" "-8.10%"
Set-TextValue $ws "G45" "11"

Set-TextValue $ws "D46" "0.00006238"
Set-TextValue $ws "E46" "-3.25%"
Set-TextValue $ws "G46" "11"

Set-TextValue $ws "E47" "-0.38%"
Set-TextValue $ws "G47" "11"

Set-TextValue $ws "D48" "0.04771"
Set-TextValue $ws "E48" "-81.05%"
Set-TextValue $ws "G48" "11"

Set-TextValue $ws "D49" "0.004201"
Set-TextValue $ws "E49" "-0.38%"
Set-TextValue $ws "G49" "11"

Set-TextValue $ws "D50" "0.00002101"
Set-TextValue $ws "E50" "-0.38%"
Set-TextValue $ws "G50" "11"

Set-TextValue $ws "D51" "0.0002001"
Set-TextValue $ws "E51" "-0.38%"
Set-TextValue $ws "G51" "11"

